$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BK (shifts old BK/BL -> BL/BM, i.e. "nom" and
# "url_produit" move one column to the right and a fresh price-history
# column is created in their place).
$ws.Columns("BK:BK").Insert()

# Header for the newly inserted column: the latest snapshot timestamp.
$ws.Range("BK1").Value2 = "2026-01-30 13:46:21"

# For every product row, the new snapshot simply repeats the last known
# price (the value that was already sitting in BJ, the previous latest
# snapshot column).
for ($r = 2; $r -le 206; $r++) {
    $prev = $ws.Cells.Item($r, 62).Value2   # column BJ = 62
    if ($prev -ne $null -and $prev -ne "") {
        $ws.Cells.Item($r, 63).Value2 = $prev   # column BK = 63
    }
}
